$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The A column values for the two data rows were previously stored as
# serial-date numbers (42809.5 / 42993.5, formatted as dates). The final
# upload instead stores the plain year "2017" as text in those cells.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2017"

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2017"
